$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Amount" value for the last expense row
$ws.Range("F7").Value = 20

# Total the amounts with a SUM formula
$ws.Range("F8").Formula = "=SUM(F4:F7)"

# Update the active cell selection to match the saved workbook state
$ws.Range("F9").Select()
